$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.849.45'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +6.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.007.79'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.61%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.48'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.08'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +13.55%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.519'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.003.43'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.84'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.157'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +7.82%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +8.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.79'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +8.29%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.866.82'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +6.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.507.93'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.98'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +7.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.006.80'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '457.84'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +6.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.96'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +8.13%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +5.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.37'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +7.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.52'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +13.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.37'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.64'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.22%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +16.86%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +15.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.61'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.85%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -6.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.00'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.36%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.993'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.81'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.73%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +12.17%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.78'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.96%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +14.62%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.64'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.48'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '389.24'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +12.38%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.796.97'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.94'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.41%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.96'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +10.98%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.06%  '
